# Update gh-pages to output generated at 456a3b4
#
# The upstream scraper re-ran: the "合肥·首届AT次元时代动漫游戏嘉年华" listing
# (2024-10-01) has expired and dropped off both the "展览" (sheet 1) and the
# "全部类型" (sheet 4) tables, so every following row shifts up by one and the
# leading index column (A) gets renumbered. A handful of "想去人数" (interest
# count) numbers also ticked up across sheets 1, 2 and 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Drop the expired first listing; everything below shifts up one row and the
# sheet dimension shrinks from A1:I9 to A1:I8 automatically.
$ws1.Rows.Item(2).Delete()

# Renumber the index column (A) for the remaining 7 data rows: 1..7
for ($r = 2; $r -le 8; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 1
}

# Refreshed "想去人数" counts
$ws1.Cells.Item(2, 6).Value = 575   # 合肥·Holic动漫游戏展
$ws1.Cells.Item(4, 6).Value = 27    # 合肥·星月动漫游戏展
$ws1.Cells.Item(5, 6).Value = 75    # 合肥·首届火影忍者同人only
$ws1.Cells.Item(7, 6).Value = 1484  # 合肥·第九届环形宇宙动漫游戏嘉年华

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances) - no rows added/removed, just two count bumps
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Cells.Item(2, 6).Value = 101   # 合肥·《四月是你的谎言》…
$ws2.Cells.Item(4, 6).Value = 1     # 合肥·豫章D乐团-《蓬莱乐，万物生》…

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - same expired listing removed as sheet 1
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Rows.Item(2).Delete()

for ($r = 2; $r -le 13; $r++) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
}

$ws4.Cells.Item(2, 6).Value = 575    # 合肥·Holic动漫游戏展
$ws4.Cells.Item(4, 6).Value = 27     # 合肥·星月动漫游戏展
$ws4.Cells.Item(5, 6).Value = 75     # 合肥·首届火影忍者同人only
$ws4.Cells.Item(7, 6).Value = 101    # 合肥·《四月是你的谎言》…
$ws4.Cells.Item(9, 6).Value = 1      # 合肥·豫章D乐团-《蓬莱乐，万物生》…
$ws4.Cells.Item(11, 6).Value = 1484  # 合肥·第九届环形宇宙动漫游戏嘉年华
